# Insert a new "posts" worksheet right after the first sheet, populate it
# with the sample CMS-style content, and make it the active/selected sheet
# (mirrors the original author re-saving with "posts" as the foreground tab
# and the old first sheet no longer tab-selected).

$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)
$posts = $wb.Worksheets.Add($null, $firstSheet)
$posts.Name = "posts"

$lorem = "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum."

$posts.Range("B1").Value = "Subhead"
$posts.Range("A1").Value = "Heading"
$posts.Range("C1").Value = "Keywords"
$posts.Range("D1").Value = "Content"

$posts.Range("A2").Value = "Getting started"
$posts.Range("C2").Value = "Tutorial"
$posts.Range("D2").Value = $lorem

$posts.Range("A3").Value = "Internediate"
$posts.Range("B3").Value = "More information"
$posts.Range("C3").Value = "Part 2"

$posts.Range("C4").Value = "Part 3"
$posts.Range("A4").Value = "Conclusion"
$posts.Range("B4").Value = "Wrapping up"

$posts.Range("D3").Value = $lorem
$posts.Range("D4").Value = $lorem

$posts.Range("B4").Select()
$posts.Activate()
